# June reporting update: add a missing Moody's rating-score row ("B3 *-" = 6)
# just above the existing "Caa1" row, shifting everything below it down by
# one row (same as Excel's Insert > Entire Row > Shift Cells Down).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 53 (pushes old row 53.. down to 54..)
$ws.Rows(53).Insert()

# Populate the newly inserted row with the new Moody's rating score
$ws.Range("A53").Value = "Moody's"
$ws.Range("B53").Value = "B3 *-"
$ws.Range("C53").Value = 6

# Match the author's final view state: scrolled so row 31 is at the top,
# with C53 as the active/selected cell.
$ws.Range("C53").Select()
$excel.ActiveWindow.ScrollRow = 31
